# WORK LOG ON 11/14/2023
# Fill in the task log for 11/14/2023 (rows 15-24 of Sheet1) and a couple
# of supporting numeric cells, then delete the stray value in B39.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- small numeric tweaks -------------------------------------------------
$ws.Range("D1").Value = 2
$ws.Range("D13").Value = 9

# --- Task Description column (B) for 11/14/2023 --------------------------
$ws.Range("B15").Value = "think about first project"
$ws.Range("B16").Value = "think about first project"
$ws.Range("B17").Value = "think about first project/lunch"
$ws.Range("B18").Value = "LED Ring schematic design"
$ws.Range("B19").Value = "LED Ring schematic design"
$ws.Range("B20").Value = "LED Ring schematic design"
$ws.Range("B21").Value = "LED Ring schematic design"
$ws.Range("B22").Value = "LED Ring schematic design"
$ws.Range("B23").Value = "search for customers"
$ws.Range("B24").Value = "search for customers"

# --- Type column (C) for 11/14/2023 ---------------------------------------
$ws.Range("C15").Value = "Support"
$ws.Range("C16").Value = "Support"
$ws.Range("C17").Value = "Support"
$ws.Range("C18").Value = "design"
$ws.Range("C19").Value = "design"
$ws.Range("C20").Value = "design"
$ws.Range("C21").Value = "design"
$ws.Range("C22").Value = "design"
$ws.Range("C23").Value = "Support"
$ws.Range("C24").Value = "Support"

# --- stray value removed from B39 -----------------------------------------
$ws.Range("B39").Clear()

# --- restore view/selection state -----------------------------------------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G14").Select() | Out-Null
